$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "GB Bin Order" column (old column F)
$ws.Columns.Item(6).Delete()

# Insert a new row for TR48 before the old TR72 row (row 29)
$ws.Rows.Item(29).Insert()

# Write out all data rows (header + data) to match the new layout
$ws.Cells.Item(1,1).Value = "Bin Label"
$ws.Cells.Item(1,2).Value = "Bin Category"
$ws.Cells.Item(1,3).Value = "Total Bins"
$ws.Cells.Item(1,4).Value = "Filled Amount"
$ws.Cells.Item(1,5).Value = "Bin Order"
$ws.Cells.Item(1,6).Value = "Bin Location"
$ws.Cells.Item(1,7).Value = "Availiability Flag"

$ws.Cells.Item(2,1).Value = "D362406"
$ws.Cells.Item(2,2).Value = "Drawer"
$ws.Cells.Item(2,3).Value = 269
$ws.Cells.Item(2,4).Value = 4.548
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = "None"
$ws.Cells.Item(2,7).Value = "Yes"

$ws.Cells.Item(3,1).Value = "D482406"
$ws.Cells.Item(3,2).Value = "Drawer"
$ws.Cells.Item(3,3).Value = 393
$ws.Cells.Item(3,4).Value = 297.9229999999996
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = "None"
$ws.Cells.Item(3,7).Value = "Yes"

$ws.Cells.Item(4,1).Value = "C361215"
$ws.Cells.Item(4,2).Value = "Clip"
$ws.Cells.Item(4,3).Value = 302
$ws.Cells.Item(4,4).Value = 302
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = "None"
$ws.Cells.Item(4,7).Value = "No"

$ws.Cells.Item(5,1).Value = "C361815"
$ws.Cells.Item(5,2).Value = "Clip"
$ws.Cells.Item(5,3).Value = 70
$ws.Cells.Item(5,4).Value = 70.00000000000001
$ws.Cells.Item(5,5).Value = 4
$ws.Cells.Item(5,6).Value = "None"
$ws.Cells.Item(5,7).Value = "No"

$ws.Cells.Item(6,1).Value = "C362415"
$ws.Cells.Item(6,2).Value = "Clip"
$ws.Cells.Item(6,3).Value = 211
$ws.Cells.Item(6,4).Value = 211
$ws.Cells.Item(6,5).Value = 5
$ws.Cells.Item(6,6).Value = "None"
$ws.Cells.Item(6,7).Value = "No"

$ws.Cells.Item(7,1).Value = "C481215"
$ws.Cells.Item(7,2).Value = "Clip"
$ws.Cells.Item(7,3).Value = 112
$ws.Cells.Item(7,4).Value = 112
$ws.Cells.Item(7,5).Value = 6
$ws.Cells.Item(7,6).Value = "None"
$ws.Cells.Item(7,7).Value = "No"

$ws.Cells.Item(8,1).Value = "C481815"
$ws.Cells.Item(8,2).Value = "Clip"
$ws.Cells.Item(8,3).Value = 186
$ws.Cells.Item(8,4).Value = 185.9999999999999
$ws.Cells.Item(8,5).Value = 7
$ws.Cells.Item(8,6).Value = "None"
$ws.Cells.Item(8,7).Value = "No"

$ws.Cells.Item(9,1).Value = "C482415"
$ws.Cells.Item(9,2).Value = "Clip"
$ws.Cells.Item(9,3).Value = 176
$ws.Cells.Item(9,4).Value = 130.2160000000002
$ws.Cells.Item(9,5).Value = 8
$ws.Cells.Item(9,6).Value = "None"
$ws.Cells.Item(9,7).Value = "Yes"

$ws.Cells.Item(10,1).Value = "B482448"
$ws.Cells.Item(10,2).Value = "Bulk"
$ws.Cells.Item(10,3).Value = 181
$ws.Cells.Item(10,4).Value = 134.3900000000055
$ws.Cells.Item(10,5).Value = 9
$ws.Cells.Item(10,6).Value = "None"
$ws.Cells.Item(10,7).Value = "Yes"

$ws.Cells.Item(11,1).Value = "B483648"
$ws.Cells.Item(11,2).Value = "Bulk"
$ws.Cells.Item(11,3).Value = 206
$ws.Cells.Item(11,4).Value = 86.49100000000027
$ws.Cells.Item(11,5).Value = 10
$ws.Cells.Item(11,6).Value = "None"
$ws.Cells.Item(11,7).Value = "Yes"

$ws.Cells.Item(12,1).Value = "B484248"
$ws.Cells.Item(12,2).Value = "Bulk"
$ws.Cells.Item(12,3).Value = 55
$ws.Cells.Item(12,4).Value = 54.99500000000015
$ws.Cells.Item(12,5).Value = 11
$ws.Cells.Item(12,6).Value = "None"
$ws.Cells.Item(12,7).Value = "No"

$ws.Cells.Item(13,1).Value = "B484848"
$ws.Cells.Item(13,2).Value = "Bulk"
$ws.Cells.Item(13,3).Value = 212
$ws.Cells.Item(13,4).Value = 52.99599999999997
$ws.Cells.Item(13,5).Value = 12
$ws.Cells.Item(13,6).Value = "None"
$ws.Cells.Item(13,7).Value = "Yes"

$ws.Cells.Item(14,1).Value = "B487248"
$ws.Cells.Item(14,2).Value = "Bulk"
$ws.Cells.Item(14,3).Value = 53
$ws.Cells.Item(14,4).Value = 52.99900000000009
$ws.Cells.Item(14,5).Value = 13
$ws.Cells.Item(14,6).Value = "None"
$ws.Cells.Item(14,7).Value = "No"

$ws.Cells.Item(15,1).Value = "B489648"
$ws.Cells.Item(15,2).Value = "Bulk"
$ws.Cells.Item(15,3).Value = 208
$ws.Cells.Item(15,4).Value = 57.87600000000001
$ws.Cells.Item(15,5).Value = 14
$ws.Cells.Item(15,6).Value = "None"
$ws.Cells.Item(15,7).Value = "Yes"

$ws.Cells.Item(16,1).Value = "B722448"
$ws.Cells.Item(16,2).Value = "Bulk"
$ws.Cells.Item(16,3).Value = 382
$ws.Cells.Item(16,4).Value = 7.804000000000002
$ws.Cells.Item(16,5).Value = 15
$ws.Cells.Item(16,6).Value = "None"
$ws.Cells.Item(16,7).Value = "Yes"

$ws.Cells.Item(17,1).Value = "B723648"
$ws.Cells.Item(17,2).Value = "Bulk"
$ws.Cells.Item(17,3).Value = 127
$ws.Cells.Item(17,4).Value = 2.945
$ws.Cells.Item(17,5).Value = 16
$ws.Cells.Item(17,6).Value = "None"
$ws.Cells.Item(17,7).Value = "Yes"

$ws.Cells.Item(18,1).Value = "B724248"
$ws.Cells.Item(18,2).Value = "Bulk"
$ws.Cells.Item(18,3).Value = 113
$ws.Cells.Item(18,4).Value = 3.921
$ws.Cells.Item(18,5).Value = 17
$ws.Cells.Item(18,6).Value = "None"
$ws.Cells.Item(18,7).Value = "Yes"

$ws.Cells.Item(19,1).Value = "B724848"
$ws.Cells.Item(19,2).Value = "Bulk"
$ws.Cells.Item(19,3).Value = 64
$ws.Cells.Item(19,4).Value = 2.07
$ws.Cells.Item(19,5).Value = 18
$ws.Cells.Item(19,6).Value = "None"
$ws.Cells.Item(19,7).Value = "Yes"

$ws.Cells.Item(20,1).Value = "B727248"
$ws.Cells.Item(20,2).Value = "Bulk"
$ws.Cells.Item(20,3).Value = 181
$ws.Cells.Item(20,4).Value = 21.767
$ws.Cells.Item(20,5).Value = 19
$ws.Cells.Item(20,6).Value = "None"
$ws.Cells.Item(20,7).Value = "Yes"

$ws.Cells.Item(21,1).Value = "B729648"
$ws.Cells.Item(21,2).Value = "Bulk"
$ws.Cells.Item(21,3).Value = 186
$ws.Cells.Item(21,4).Value = 5.990999999999999
$ws.Cells.Item(21,5).Value = 20
$ws.Cells.Item(21,6).Value = "None"
$ws.Cells.Item(21,7).Value = "Yes"

$ws.Cells.Item(22,1).Value = "B962448"
$ws.Cells.Item(22,2).Value = "Bulk"
$ws.Cells.Item(22,3).Value = 178
$ws.Cells.Item(22,4).Value = 8.300999999999998
$ws.Cells.Item(22,5).Value = 21
$ws.Cells.Item(22,6).Value = "None"
$ws.Cells.Item(22,7).Value = "Yes"

$ws.Cells.Item(23,1).Value = "B963648"
$ws.Cells.Item(23,2).Value = "Bulk"
$ws.Cells.Item(23,3).Value = 323
$ws.Cells.Item(23,4).Value = 1.766
$ws.Cells.Item(23,5).Value = 22
$ws.Cells.Item(23,6).Value = "None"
$ws.Cells.Item(23,7).Value = "Yes"

$ws.Cells.Item(24,1).Value = "B964248"
$ws.Cells.Item(24,2).Value = "Bulk"
$ws.Cells.Item(24,3).Value = 378
$ws.Cells.Item(24,4).Value = 0.753
$ws.Cells.Item(24,5).Value = 23
$ws.Cells.Item(24,6).Value = "None"
$ws.Cells.Item(24,7).Value = "Yes"

$ws.Cells.Item(25,1).Value = "B964848"
$ws.Cells.Item(25,2).Value = "Bulk"
$ws.Cells.Item(25,3).Value = 161
$ws.Cells.Item(25,4).Value = 0.744
$ws.Cells.Item(25,5).Value = 24
$ws.Cells.Item(25,6).Value = "None"
$ws.Cells.Item(25,7).Value = "Yes"

$ws.Cells.Item(26,1).Value = "B967248"
$ws.Cells.Item(26,2).Value = "Bulk"
$ws.Cells.Item(26,3).Value = 214
$ws.Cells.Item(26,4).Value = 0.183
$ws.Cells.Item(26,5).Value = 25
$ws.Cells.Item(26,6).Value = "None"
$ws.Cells.Item(26,7).Value = "Yes"

$ws.Cells.Item(27,1).Value = "B969648"
$ws.Cells.Item(27,2).Value = "Bulk"
$ws.Cells.Item(27,3).Value = 121
$ws.Cells.Item(27,4).Value = 0.248
$ws.Cells.Item(27,5).Value = 26
$ws.Cells.Item(27,6).Value = "None"
$ws.Cells.Item(27,7).Value = "Yes"

$ws.Cells.Item(28,1).Value = "BR484816"
$ws.Cells.Item(28,2).Value = "Battery"
$ws.Cells.Item(28,3).Value = 233
$ws.Cells.Item(28,4).Value = 6.613899999999998
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(28,6).Value = "None"
$ws.Cells.Item(28,7).Value = "Yes"

$ws.Cells.Item(29,1).Value = "TR48"
$ws.Cells.Item(29,2).Value = "Tire"
$ws.Cells.Item(29,3).Value = 100
$ws.Cells.Item(29,4).Value = 0
$ws.Cells.Item(29,5).Value = 0
$ws.Cells.Item(29,6).Value = "None"
$ws.Cells.Item(29,7).Value = "Yes"

$ws.Cells.Item(30,1).Value = "TR72"
$ws.Cells.Item(30,2).Value = "Tire"
$ws.Cells.Item(30,3).Value = 108
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(30,6).Value = "None"
$ws.Cells.Item(30,7).Value = "Yes"

$ws.Cells.Item(31,1).Value = "BC967248"
$ws.Cells.Item(31,2).Value = "Bumper Cover"
$ws.Cells.Item(31,3).Value = 275
$ws.Cells.Item(31,4).Value = 27.73579999999998
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(31,6).Value = "None"
$ws.Cells.Item(31,7).Value = "Yes"

$ws.Cells.Item(32,1).Value = "BH967280"
$ws.Cells.Item(32,2).Value = "Hood"
$ws.Cells.Item(32,3).Value = 295
$ws.Cells.Item(32,4).Value = 2.8585
$ws.Cells.Item(32,5).Value = 0
$ws.Cells.Item(32,6).Value = "None"
$ws.Cells.Item(32,7).Value = "Yes"

$ws.Cells.Item(33,1).Value = "HS06"
$ws.Cells.Item(33,2).Value = "Hanging"
$ws.Cells.Item(33,3).Value = 377
$ws.Cells.Item(33,4).Value = 377
$ws.Cells.Item(33,5).Value = 0
$ws.Cells.Item(33,6).Value = "None"
$ws.Cells.Item(33,7).Value = "No"

$ws.Cells.Item(34,1).Value = "HS12"
$ws.Cells.Item(34,2).Value = "Hanging"
$ws.Cells.Item(34,3).Value = 364
$ws.Cells.Item(34,4).Value = 364
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(34,6).Value = "None"
$ws.Cells.Item(34,7).Value = "No"
